# Deloitte.xlsx (re-)upload: refresh the logged "Request"/"Response" JSON
# payloads on Sheet3 with a newer reqres.in API automation run (Janet
# Weaver fixture) and move the sheet selection off the stale "B10" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / col A ("Request" body) - pretty-printed JSON payload that was
# sent to the API.
$requestBody = '{' + "`n" + '"email":"janet.weaver@reqres.in",' + "`n" + '"first_name":"Janet",' + "`n" + '"last_name":"Weaver"' + "`n" + '}'
$ws.Range("A2").Value = $requestBody

# Row 2 / col B ("Response" body) - compact JSON payload returned by the
# API for the latest run in the series.
$responseBody = '{"email":"janet.weaver@reqres.in","first_name":"Janet","last_name":"Weaver","id":"973","createdAt":"2024-07-17T18:39:23.581Z"}'
$ws.Range("B2").Value = $responseBody

# The longer pretty-printed request body now wraps across more lines, so
# the row grows from 52.8pt to 79.2pt.
$ws.Rows.Item(2).RowHeight = 79.2

# Selection moves from the old B10 cell to A2.
$ws.Range("A2").Select() | Out-Null
